$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.326.93"
$ws.Range("E2").Value = "  -2.10%  "

$ws.Range("D3").Value = "2.364.02"
$ws.Range("E3").Value = "  -0.40%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'328.12"
$ws.Range("E5").Value = "  +4.66%  "

$ws.Range("D6").Value = "'99.85"
$ws.Range("E6").Value = "  -6.77%  "

$ws.Range("D7").Value = "'0.634"
$ws.Range("E7").Value = "  -0.84%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.616"
$ws.Range("E9").Value = "  -1.85%  "

$ws.Range("D10").Value = "'40.01"
$ws.Range("E10").Value = "  -7.05%  "

$ws.Range("D11").Value = "'0.0919"
$ws.Range("E11").Value = "  -2.18%  "

$ws.Range("D12").Value = "'8.44"
$ws.Range("E12").Value = "  -5.38%  "

$ws.Range("D13").Value = "'1.00"
$ws.Range("E13").Value = "  -5.45%  "

$ws.Range("D14").Value = "'0.105"
$ws.Range("E14").Value = "  +0.13%  "

$ws.Range("D15").Value = "'16.22"
$ws.Range("E15").Value = "  -2.24%  "

$ws.Range("D16").Value = "2.719.26"
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("D17").Value = "2.365.34"
$ws.Range("E17").Value = "  -0.47%  "

$ws.Range("D18").Value = "42.413.95"
$ws.Range("E18").Value = "  -1.85%  "

$ws.Range("D19").Value = "'7.75"
$ws.Range("E19").Value = "  +5.42%  "

$ws.Range("E20").Value = "  -2.32%  "

$ws.Range("D21").Value = "'3.73"
$ws.Range("E21").Value = "  +7.35%  "

$ws.Range("D22").Value = "'74.73"
$ws.Range("E22").Value = "  -1.18%  "

$ws.Range("D23").Value = "'276.33"
$ws.Range("E23").Value = "  +8.88%  "

$ws.Range("D24").Value = "'2.28"
$ws.Range("E24").Value = "  -9.35%  "

$ws.Range("D25").Value = "'9.58"
$ws.Range("E25").Value = "  +7.52%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").Value = "'11.41"
$ws.Range("E27").Value = "  -5.45%  "

$ws.Range("D28").Value = "'23.62"
$ws.Range("E28").Value = "  +3.20%  "

$ws.Range("E29").Value = "  -1.48%  "

$ws.Range("D30").Value = "'174.12"
$ws.Range("E30").Value = "  +0.58%  "

$ws.Range("E31").Value = "  -2.22%  "

$ws.Range("D32").Value = "'0.0897"
$ws.Range("E32").Value = "  -1.46%  "

$ws.Range("D33").Value = "'35.08"
$ws.Range("E33").Value = "  -10.26%  "

$ws.Range("D34").Value = "'6.00"
$ws.Range("E34").Value = "  +2.46%  "

$ws.Range("E35").Value = "  +1.80%  "

$ws.Range("D36").Value = "'4.52"
$ws.Range("E36").Value = "  -9.46%  "

$ws.Range("D37").Value = "'0.0357"
$ws.Range("E37").Value = "  -5.29%  "

$ws.Range("D38").Value = "'2.92"
$ws.Range("E38").Value = "  +4.56%  "

$ws.Range("D39").Value = "'3.84"
$ws.Range("E39").Value = "  -6.31%  "

$ws.Range("E40").Value = "  +0.91%  "

$ws.Range("D41").Value = "'1.52"
$ws.Range("E41").Value = "  -1.30%  "

$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "'68.94"
$ws.Range("E42").Value = "  -4.72%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.227"
$ws.Range("E43").Value = "  -2.51%  "

$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D45").Value = "'90.79"
$ws.Range("E45").Value = "  +30.50%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'115.83"
$ws.Range("E46").Value = "  +2.50%  "

$ws.Range("D47").Value = "'11.90"
$ws.Range("E47").Value = "  -4.04%  "

$ws.Range("D48").Value = "'5.45"
$ws.Range("E48").Value = "  -4.12%  "

$ws.Range("D49").Value = "'9.08"
$ws.Range("E49").Value = "  -2.67%  "

$ws.Range("D50").Value = "1.594.19"
$ws.Range("E50").Value = "  +6.37%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.100"
$ws.Range("E51").Value = "  +0.62%  "
